$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 237, shifting existing rows 237-290 down to 238-291
$ws.Rows.Item(237).Insert()

# Populate the newly inserted row 237 with the new record's data
$r = 237
$ws.Cells.Item($r, 1).Value = 4
$ws.Cells.Item($r, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($r, 3).Value = "Los Lagos"
$ws.Cells.Item($r, 4).Value = 44798
$ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 5).Value = 10
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108005
$ws.Cells.Item($r, 10).Value = "Piña"
$ws.Cells.Item($r, 11).Value = "Caramelo"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 100
$ws.Cells.Item($r, 14).Value = 23000
$ws.Cells.Item($r, 15).Value = 23000
$ws.Cells.Item($r, 16).Value = 23000
$ws.Cells.Item($r, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item($r, 18).Value = "Ecuador"
$ws.Cells.Item($r, 19).Value = 1917
$ws.Cells.Item($r, 20).Value = 12
